# always keep family infos, keep orig names also if not changed by GBIF
#
# GBIF sometimes can only resolve a name down to genus or family level
# (column A = "Species GBIF"), while the original observation (column B =
# "Species Original") is more specific (e.g. "Acer" vs "Acer sp."). In
# those cases, append " species" to the GBIF name so the higher-rank
# match is clearly labelled (e.g. "Acer" -> "Acer species").
#
# Rule: a row needs the suffix when the GBIF name (column A) is a single
# taxon word (genus/family only, no species epithet) AND it differs from
# the original name (column B) - i.e. GBIF did not resolve the name to
# species level, so we keep the family/genus info but mark it explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()

    if ($null -eq $a) { continue }

    $aTrim = $a.Trim()
    if ($aTrim.Length -eq 0) { continue }

    $isSingleWord = -not $aTrim.Contains(" ")
    $differsFromOriginal = $aTrim -ne $b

    if ($isSingleWord -and $differsFromOriginal) {
        $ws.Cells.Item($r, 1).Value = "$aTrim species"
    }
}
